# Apply the cryptos-list price/volume/coin updates described in the commit.
# Every touched cell is forced to text format ("@") before the value is written so
# that numeric-looking strings (e.g. "7.00", "0.812", "19.39") are stored verbatim as
# text, exactly like the original inline-string cells, instead of being coerced to
# numbers (which would drop trailing zeros / introduce floating-point noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.417.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.310.86"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.32"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.36%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.97%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.668.35"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.317.14"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.320.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.61"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.63"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.43%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.94%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.94"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.63"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.24%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.62"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.64%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.48%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.17%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.39"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.976.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0289"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.99%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.15%  "
